$d = $word.ActiveDocument

# The word "gewisse" (in "... op een gewisse dood. ") was originally split
# across three runs:  " " + "ge" + "wisse dood. "
# The edit merges those three runs into a single run containing
# " wisse dood. " (the leading "ge" is dropped, turning "gewisse" into
# "wisse"). "gewisse" is unique in the document, so searching/replacing on
# it precisely collapses just those three runs into one, without touching
# the preceding run ("... op de berg te wachten op een").
$d.Content.Find.Execute("gewisse dood.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "wisse dood.", 2)
